$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 (shifts existing rows 33:95 down to 34:96)
$ws.Rows.Item(33).Insert()

# Copy the repeating/constant field values from the row above (now row 34,
# originally row 33) into the freshly inserted row 33.
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = "Femacal de La Calera"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 100112052
$ws.Range("G33").Value = "Albahaca"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("N33").Value = "$/docena de matas"
$ws.Range("O33").Value = "Provincia de Quillota"
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = "Hortaliza"

# New weekly record values
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D33").Value = 44540
$ws.Range("J33").Value = 105
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 4500
$ws.Range("M33").Value = 4262
$ws.Range("P33").Value = 710
